# Re-crawl update for the 2022-09-02 Coop "toilet paper" listing:
#
#   1. The crawler re-queried the page later the same day, so every row's
#      `timestamp` (column O) moves from 07:07:59 to 21:00:15.
#   2. Between the two crawls, three neighbouring product pairs swapped
#      positions in the listing (rows 18/19, 25/26, 28/29) - their A:N
#      content (id, title, href, quantity, ratings, brand, price, ...)
#      is exchanged, row-for-row.
#
# Columns E (ratingAmount) and F (ratingValue) are numeric; everything else
# in A:N is text (ids/prices/etc. are stored as text in the source data, not
# numbers - e.g. the price "17.70" must stay "17.70", not become 17.7).
# Reading through `.Text` keeps the original formatted representation, and
# writing text back with a leading apostrophe stops Excel from re-coercing
# numeric-looking strings (ids, prices, quantities) into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericCols = @(5, 6)   # E = ratingAmount, F = ratingValue
$lastTextCol = 14        # A..N

function Get-RowText {
    param($row)
    $vals = @()
    for ($c = 1; $c -le $lastTextCol; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Text
    }
    return $vals
}

function Set-RowFromText {
    param($row, $vals)
    for ($c = 1; $c -le $lastTextCol; $c++) {
        $v = $vals[$c - 1]
        if ($numericCols -contains $c) {
            if ([string]::IsNullOrEmpty($v)) {
                # Blank rating cells are empty *text* in the source, not 0.
                $ws.Cells.Item($row, $c).Value = "'"
            } else {
                $ws.Cells.Item($row, $c).Value = [double]$v
            }
        } else {
            $ws.Cells.Item($row, $c).Value = "'" + $v
        }
    }
}

function Swap-Rows {
    param($r1, $r2)
    $v1 = Get-RowText $r1
    $v2 = Get-RowText $r2
    Set-RowFromText $r1 $v2
    Set-RowFromText $r2 $v1
}

Swap-Rows 18 19
Swap-Rows 25 26
Swap-Rows 28 29

# --- refresh the crawl timestamp on every data row (2-31) ------------------
$newTimestamp = "2022-09-02 21:00:15"
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 15).Value = "'" + $newTimestamp
}

Write-Host "done"
